# Update the form date on A2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Date: 08 Nov 2025"

# Append a new data row (row 11) mirroring the existing rows' layout/style
$ws.Range("A10:F10").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = "*dummy*"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "Relocating to CSF"
$ws.Range("F11").Value = ""

$ws.Rows.Item(11).RowHeight = 13
